$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.435.71'
$ws.Range("E2").Value = '  +1.53%  '

$ws.Range("D3").Value = '2.289.74'
$ws.Range("E3").Value = '  +3.06%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '496.50'
$ws.Range("E5").Value = '  +2.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.71'
$ws.Range("E6").Value = '  +2.83%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.531'
$ws.Range("E8").Value = '  +2.83%  '

$ws.Range("D9").Value = '2.287.48'
$ws.Range("E9").Value = '  +2.56%  '

$ws.Range("E10").Value = '  +4.43%  '

$ws.Range("E11").Value = '  +2.56%  '

$ws.Range("E12").Value = '  +4.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.65'
$ws.Range("E13").Value = '  -0.26%  '

$ws.Range("D14").Value = '2.693.35'
$ws.Range("E14").Value = '  +3.02%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.83'
$ws.Range("E15").Value = '  +4.13%  '

$ws.Range("D16").Value = '54.339.40'
$ws.Range("E16").Value = '  +1.49%  '

$ws.Range("E17").Value = '  +1.80%  '

$ws.Range("D18").Value = '2.289.60'
$ws.Range("E18").Value = '  +3.47%  '

$ws.Range("E19").Value = '  +5.49%  '

$ws.Range("E20").Value = '  +4.42%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '301.86'
$ws.Range("E21").Value = '  +1.36%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.45'
$ws.Range("E22").Value = '  +5.95%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("E24").Value = '  -1.86%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.63'
$ws.Range("E25").Value = '  -0.99%  '

$ws.Range("E26").Value = '  +1.57%  '

$ws.Range("E27").Value = '  +3.08%  '

$ws.Range("D28").Value = '2.390.58'
$ws.Range("E28").Value = '  +3.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.149'
$ws.Range("E29").Value = '  +4.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.09'
$ws.Range("E30").Value = '  +1.86%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.81'
$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D32").Value = '0.0₃0692'
$ws.Range("E32").Value = '  +2.75%  '

$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.61'
$ws.Range("E33").Value = '  +2.19%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.90'
$ws.Range("E34").Value = '  +3.36%  '

$ws.Range("E35").Value = '  +0.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.16%  '

$ws.Range("E37").Value = '  +2.53%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.73'
$ws.Range("E38").Value = '  +2.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.913'
$ws.Range("E39").Value = '  +10.08%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.20'
$ws.Range("E40").Value = '  +4.48%  '

$ws.Range("E41").Value = '  +4.36%  '

$ws.Range("E42").Value = '  -0.64%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.41'
$ws.Range("E43").Value = '  +3.62%  '

$ws.Range("E44").Value = '  +2.85%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.37'
$ws.Range("E45").Value = '  +3.42%  '

$ws.Range("E46").Value = '  +3.92%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.80'
$ws.Range("E47").Value = '  +3.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0889'
$ws.Range("E48").Value = '  +1.53%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.547'
$ws.Range("E49").Value = '  +2.95%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '239.22'
$ws.Range("E50").Value = '  +4.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0485'
$ws.Range("E51").Value = '  +3.80%  '
